$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 630
$ws1.Range("F4").Value = 646
$ws1.Range("F5").Value = 565
$ws1.Range("F7").Value = 2783
$ws1.Range("F9").Value = 7823
$ws1.Range("F11").Value = 470
$ws1.Range("F12").Value = 43
$ws1.Range("F13").Value = 349

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 630
$ws4.Range("F4").Value = 646
$ws4.Range("F5").Value = 565
$ws4.Range("F9").Value = 2783
$ws4.Range("F11").Value = 7823
$ws4.Range("F13").Value = 470
$ws4.Range("F14").Value = 43
$ws4.Range("F17").Value = 349
